$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update shared string "Resolving-Mac" -> "Neutrophils" (Target cluster, D5 and D9)
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("D9").Value = "Neutrophils"

# Update recalculated TPM-derived numeric values
$ws.Range("H2").Value = 4.501409000000001
$ws.Range("M2").Value = 2.820332666666667
$ws.Range("N2").Value = 8.460998
$ws.Range("O2").Value = 0.5374031335545234
$ws.Range("P2").Value = 0.5374031335545234
$ws.Range("Q2").Value = 4.231823616242445
$ws.Range("R2").Value = 38.086412546182
$ws.Range("S2").Value = 0.5343852766749654
$ws.Range("T2").Value = 0.5343852766749655

$ws.Range("H3").Value = 4.501409000000001
$ws.Range("O3").Value = 0.2710731659200333
$ws.Range("P3").Value = 0.2710731659200333
$ws.Range("S3").Value = 0.2695509194581928
$ws.Range("T3").Value = 0.2695509194581929

$ws.Range("H4").Value = 4.501409000000001
$ws.Range("M4").Value = 0.8926926666666667
$ws.Range("N4").Value = 2.678078
$ws.Range("O4").Value = 0.170099024855393
$ws.Range("P4").Value = 0.170099024855393
$ws.Range("Q4").Value = 1.339458267989111
$ws.Range("R4").Value = 12.055124411902
$ws.Range("S4").Value = 0.1691438117568563
$ws.Range("T4").Value = 0.1691438117568564

$ws.Range("H5").Value = 4.501409000000001
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1124383333333333
$ws.Range("N5").Value = 0.337315
$ws.Range("O5").Value = 0.02142467567005028
$ws.Range("P5").Value = 0.02142467567005028
$ws.Range("Q5").Value = 0.1687103085372223
$ws.Range("R5").Value = 1.518392776835
$ws.Range("S5").Value = 0.02130436262975313
$ws.Range("T5").Value = 0.02130436262975313

$ws.Range("I6").Value = 0.005615629480232302
$ws.Range("J6").Value = 0.005615629480232303
$ws.Range("M6").Value = 2.820332666666667
$ws.Range("N6").Value = 8.460998
$ws.Range("O6").Value = 0.5374031335545234
$ws.Range("P6").Value = 0.5374031335545234
$ws.Range("Q6").Value = 0.02389855890644444
$ws.Range("R6").Value = 0.215087030158
$ws.Range("S6").Value = 0.003017856879557999
$ws.Range("T6").Value = 0.003017856879557999

$ws.Range("I7").Value = 0.005615629480232302
$ws.Range("J7").Value = 0.005615629480232303
$ws.Range("O7").Value = 0.2710731659200333
$ws.Range("P7").Value = 0.2710731659200333
$ws.Range("S7").Value = 0.001522246461840441
$ws.Range("T7").Value = 0.001522246461840441

$ws.Range("I8").Value = 0.005615629480232302
$ws.Range("J8").Value = 0.005615629480232303
$ws.Range("M8").Value = 0.8926926666666667
$ws.Range("N8").Value = 2.678078
$ws.Range("O8").Value = 0.170099024855393
$ws.Range("P8").Value = 0.170099024855393
$ws.Range("Q8").Value = 0.007564380093111111
$ws.Range("R8").Value = 0.068079420838
$ws.Range("S8").Value = 0.0009552130985367123
$ws.Range("T8").Value = 0.0009552130985367125

$ws.Range("I9").Value = 0.005615629480232302
$ws.Range("J9").Value = 0.005615629480232303
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1124383333333333
$ws.Range("N9").Value = 0.337315
$ws.Range("O9").Value = 0.02142467567005028
$ws.Range("P9").Value = 0.02142467567005028
$ws.Range("Q9").Value = 0.0009527649572222222
$ws.Range("R9").Value = 0.008574884615
$ws.Range("S9").Value = 0.0001203130402971501
$ws.Range("T9").Value = 0.0001203130402971501
